$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$data = @(
    @(2, 8, 8),
    @(3, 8, 9),
    @(4, 9, 9),
    @(5, 9, 9),
    @(6, 9, 9),
    @(7, 8, 9),
    @(8, 9, 9),
    @(9, 9, 9),
    @(10, 8, 9),
    @(11, 10, 11),
    @(12, 8, 9),
    @(13, 9, 9),
    @(14, 9, 9),
    @(15, 9, 9),
    @(16, 9, 9),
    @(17, 9, 9),
    @(18, 8, 9),
    @(19, 9, 9),
    @(20, 9, 9),
    @(21, 9, 9),
    @(22, 7, 8),
    @(23, 8, 9),
    @(24, 8, 8),
    @(25, 9, 9),
    @(26, 8, 8),
    @(27, 8, 8),
    @(28, 8, 8),
    @(29, 9, 9),
    @(30, 7, 8),
    @(31, 9, 9),
    @(32, 7, 7),
    @(33, 8, 8),
    @(34, 8, 8),
    @(35, 8, 8),
    @(36, 6, 7),
    @(37, 7, 7),
    @(38, 7, 8),
    @(39, 10, 10),
    @(40, 6, 6),
    @(41, 8, 8),
    @(42, 7, 7),
    @(43, 7, 7),
    @(44, 7, 7),
    @(45, 9, 9),
    @(46, 7, 7),
    @(47, 6, 7),
    @(48, 7, 8),
    @(49, 8, 8),
    @(50, 9, 9),
    @(51, 8, 8),
    @(52, 9, 9),
    @(53, 8, 9),
    @(54, 9, 9),
    @(55, 7, 7),
    @(56, 5, 6),
    @(57, 7, 8),
    @(58, 5, 5),
    @(59, 6, 7),
    @(60, 7, 8),
    @(61, 7, 7),
    @(62, 7, 7),
    @(63, 6, 6),
    @(64, 8, 8),
    @(65, 6, 7),
    @(66, 2, 3),
    @(67, 6, 6),
    @(68, 4, 4),
    @(69, 8, 8),
    @(70, 7, 7)
)

foreach ($row in $data) {
    $r = $row[0]
    $i0 = $row[1]
    $if = $row[2]
    $ws.Cells.Item($r, 9).Value = $i0
    $ws.Cells.Item($r, 10).Value = $if
}
